$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.097.53'
$ws.Range("E2").Value = '  +6.66%  '
$ws.Range("D3").Value = '3.017.96'
$ws.Range("E3").Value = '  +3.82%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''585.30'
$ws.Range("E5").Value = '  +2.90%  '
$ws.Range("D6").Value = '''162.28'
$ws.Range("E6").Value = '  +12.90%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.014.36'
$ws.Range("E8").Value = '  +3.82%  '
$ws.Range("E9").Value = '  +3.60%  '
$ws.Range("D10").Value = '''6.74'
$ws.Range("E10").Value = '  -3.11%  '
$ws.Range("E11").Value = '  +5.45%  '
$ws.Range("E12").Value = '  +6.55%  '
$ws.Range("D13").Value = '''0.0000255'
$ws.Range("E13").Value = '  +7.91%  '
$ws.Range("D14").Value = '''34.80'
$ws.Range("E14").Value = '  +7.13%  '
$ws.Range("E15").Value = '  -0.65%  '
$ws.Range("D16").Value = '66.034.01'
$ws.Range("E16").Value = '  +6.67%  '
$ws.Range("D17").Value = '3.518.01'
$ws.Range("E17").Value = '  +3.84%  '
$ws.Range("D18").Value = '''6.95'
$ws.Range("E18").Value = '  +7.00%  '
$ws.Range("D19").Value = '3.015.61'
$ws.Range("E19").Value = '  +3.98%  '
$ws.Range("D20").Value = '''457.91'
$ws.Range("E20").Value = '  +6.61%  '
$ws.Range("D21").Value = '''13.95'
$ws.Range("E21").Value = '  +7.18%  '
$ws.Range("D22").Value = '''0.690'
$ws.Range("E22").Value = '  +5.92%  '
$ws.Range("D23").Value = '''7.41'
$ws.Range("E23").Value = '  +7.89%  '
$ws.Range("E24").Value = '  +4.57%  '
$ws.Range("D25").Value = '''2.30'
$ws.Range("E25").Value = '  +13.50%  '
$ws.Range("D26").Value = '''12.42'
$ws.Range("E26").Value = '  +3.52%  '
$ws.Range("D27").Value = '''10.65'
$ws.Range("E27").Value = '  +3.92%  '
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("D29").Value = '''8.08'
$ws.Range("E29").Value = '  +16.20%  '
$ws.Range("D30").Value = '''2.36'
$ws.Range("E30").Value = '  +17.45%  '
$ws.Range("D31").Value = '''0.0000106'
$ws.Range("E31").Value = '  -6.50%  '
$ws.Range("D32").Value = '''2.61'
$ws.Range("E32").Value = '  +4.33%  '
$ws.Range("D33").Value = '''27.34'
$ws.Range("E33").Value = '  +6.77%  '
$ws.Range("D34").Value = '''0.112'
$ws.Range("E34").Value = '  +5.46%  '
$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").Value = '''0.996'
$ws.Range("E36").Value = '  +4.60%  '
$ws.Range("D37").Value = '''5.84'
$ws.Range("E37").Value = '  +8.31%  '
$ws.Range("E38").Value = '  +15.03%  '
$ws.Range("D39").Value = '''3.03'
$ws.Range("E39").Value = '  +4.54%  '
$ws.Range("E40").Value = '  +2.20%  '
$ws.Range("D41").Value = '''0.311'
$ws.Range("E41").Value = '  +16.70%  '
$ws.Range("E42").Value = '  +7.66%  '
$ws.Range("D43").Value = '''43.42'
$ws.Range("E43").Value = '  +5.55%  '
$ws.Range("D44").Value = '''8.46'
$ws.Range("E44").Value = '  +3.64%  '
$ws.Range("D45").Value = '''396.56'
$ws.Range("E45").Value = '  +13.54%  '
$ws.Range("D46").Value = '''0.0360'
$ws.Range("E46").Value = '  +7.19%  '
$ws.Range("D47").Value = '2.804.04'
$ws.Range("E47").Value = '  +3.38%  '
$ws.Range("D48").Value = '''133.53'
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("E50").Value = '  +10.82%  '
$ws.Range("E51").Value = '  +4.43%  '
